$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for existing rows 2-7 (columns A-K) and new rows 8-17.
# Row layout: A, B, C, D, E, F, G, H, I, J, K
$data = @{
    2  = @(39, 40, 354.1473291618483, 0.01111112, 0.01111112, 0.02222224, 531.2209937427724, 0.31, 0.44, 0.18, 0.32)
    3  = @(40, 41, 356.5382138814032, 0, 0, 0, 534.8073208221048, 0.33, 0.5, 0.12, 0.36)
    4  = @(41, 42, 444.5860015701382, 0, 0, 0, 666.8790023552073, 0.3, 0.45, 0.15, 0.35)
    5  = @(42, 43, 515.7790399549028, 0, 0, 0, 773.6685599323542, 0.4, 0.48, 0.19, 0.4)
    6  = @(43, 44, 517.557383961184, 0, 0, 0, 776.336075941776, 0.35, 0.43, 0.19, 0.39)
    7  = @(44, 45, 423.5225047846375, 0, 0, 0, 635.2837571769562, 0.38, 0.46, 0.11, 0.31)
    8  = @(45, 46, 330.2780007219939, 0, 0, 0, 495.4170010829909, 0.38, 0.44, 0.11, 0.35)
    9  = @(46, 47, 301.7454591101147, 0, 0, 0, 452.6181886651721, 0.34, 0.4, 0.2, 0.3)
    10 = @(47, 48, 255.9827400151617, 0, 0, 0, 383.9741100227426, 0.32, 0.42, 0.2, 0.3)
    11 = @(48, 49, 192.8515277922028, 0, 0, 0, 289.2772916883042, 0.34, 0.48, 0.18, 0.35)
    12 = @(49, 50, 138.809629379121, 0, 0, 0, 208.2144440686815, 0.4, 0.48, 0.2, 0.34)
    13 = @(50, 51, 111.1662597703822, 0, 0, 0, 166.7493896555733, 0.4, 0.5, 0.11, 0.33)
    14 = @(51, 52, 101.4841646250765, 0, 0, 0, 152.2262469376147, 0.32, 0.41, 0.13, 0.31)
    15 = @(52, 53, 95.53659189296098, 0, 0, 0, 143.3048878394415, 0.38, 0.46, 0.17, 0.33)
    16 = @(53, 54, 95.25996060309471, 0, 0, 0, 142.8899409046421, 0.38, 0.47, 0.12, 0.39)
    17 = @(54, 55, 101.009939556736, 0, 0, 0, 151.514909335104, 0.4, 0.48, 0.17, 0.31)
}

$cols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J", "K")

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $vals[$i]
    }
}

# Column A carries the bordered / centered-bold format (same style as the
# existing A2:A7 cells) for every newly added row (8-17). Copy the format
# from an already-styled cell instead of touching individual format
# properties, so no stray/unused style entries get minted.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A8:A17").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
